$p = $ppt.ActivePresentation

# --- Remove the last slide ("En devoir", slide 17) ---
$p.Slides.Item($p.Slides.Count).Delete()

# --- Refresh the cached "datetimeFigureOut" date field text from
#     2022-04-26 to 2022-04-29 across the slide master and every
#     slide layout that carries a Date placeholder. ---
$newDate = "2022-04-29"

$m = $p.SlideMaster
$m.Shapes.Item(3).TextFrame.TextRange.Text = $newDate

$dateShapeIndexByLayout = @{1=3; 2=3; 3=3; 4=4; 5=1; 6=1; 7=1}
for ($i = 1; $i -le $m.CustomLayouts.Count; $i++) {
  $layout = $m.CustomLayouts.Item($i)
  $idx = $dateShapeIndexByLayout[$i]
  if ($idx -ne $null) {
    $layout.Shapes.Item($idx).TextFrame.TextRange.Text = $newDate
  }
}
